# Angepasste Foerderfloskeln vom Schulamt Dortmund
# Replaces the Wuppertal-flavoured funding-remark text blocks (rows 2-20 in the
# original sheet) with the revised Dortmund set, extending the table down to
# row 26 and leaving column C ("Art") blank for the new Hoeren-und-Kommunikation /
# Sehen rows (9-14) which have no "FSP" marker in the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$Vorname$ wurde im Förderschwerpunkt Sprache sonderpädagogisch gefördert und im Bildungsgang der Grundschule unterrichtet.'
$ws.Range("C2").Value = 'FSP'
$ws.Range("B3").Value = '$Vorname$ wurde im Förderschwerpunkt Sprache sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Lernen unterrichtet.'
$ws.Range("C3").Value = 'FSP'
$ws.Range("B4").Value = '$Vorname$ wurde im Förderschwerpunkt Emotionale und soziale Entwicklung sonderpädagogisch gefördert und im Bildungsgang der Grundschule unterrichtet.'
$ws.Range("C4").Value = 'FSP'
$ws.Range("B5").Value = '$Vorname$ wurde im Förderschwerpunkt Emotionale und soziale Entwicklung sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Lernen unterrichtet.'
$ws.Range("C5").Value = 'FSP'
$ws.Range("B6").Value = '$Vorname$ wurde im Förderschwerpunkt Körperliche und motorische Entwicklung sonderpädagogisch gefördert und im Bildungsgang der Grundschule unterrichtet.'
$ws.Range("C6").Value = 'FSP'
$ws.Range("B7").Value = '$Vorname$ wurde im Förderschwerpunkt Körperliche und motorische Entwicklung sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Lernen unterrichtet.'
$ws.Range("C7").Value = 'FSP'
$ws.Range("B8").Value = '$Vorname$ wurde im Förderschwerpunkt Körperliche und motorische Entwicklung sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Geistige Entwicklung unterrichtet.'
$ws.Range("C8").Value = 'FSP'
$ws.Range("B9").Value = '$Vorname$ wurde im Förderschwerpunkt Hören und Kommunikation sonderpädagogisch gefördert und im Bildungsgang der Grundschule unterrichtet.'
$ws.Range("C9").ClearContents()
$ws.Range("B10").Value = '$Vorname$ wurde im Förderschwerpunkt Hören und Kommunikation sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Lernen unterrichtet.'
$ws.Range("C10").ClearContents()
$ws.Range("B11").Value = '$Vorname$ wurde im Förderschwerpunkt Hören und Kommunikation sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Geistige Entwicklung unterrichtet.'
$ws.Range("C11").ClearContents()
$ws.Range("B12").Value = '$Vorname$ wurde im Förderschwerpunkt Sehen sonderpädagogisch gefördert und im Bildungsgang der Grundschule unterrichtet.'
$ws.Range("C12").ClearContents()
$ws.Range("B13").Value = '$Vorname$ wurde im Förderschwerpunkt Sehen sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Lernen unterrichtet.'
$ws.Range("C13").ClearContents()
$ws.Range("B14").Value = '$Vorname$ wurde im Förderschwerpunkt Sehen sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Geistige Entwicklung unterrichtet.'
$ws.Range("C14").ClearContents()
$ws.Range("B15").Value = '$Vorname$ wurde im Förderschwerpunkt Lernen sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Lernen unterrichtet.'
$ws.Range("C15").Value = 'FSP'
$ws.Range("B16").Value = '$Vorname$ wurde im Förderschwerpunkt Geistige Entwicklung sonderpädagogisch gefördert und im zieldifferenten Bildungsgang Geistige Entwicklung unterrichtet.'
$ws.Range("C16").Value = 'FSP'
$ws.Range("B17").Value = '$Vorname$ hat gemäß § 18 AO-SF durch die Entscheidung des Schulamtes **Wuppertal** der Bezirksregierung **Düsseldorf** vom **tt.mm.jjjj** keinen Bedarf an sonderpädagogischer Unterstützung mehr.'
$ws.Range("C17").Value = 'FSP'
$ws.Range("B18").Value = '$Vorname$ wechselt gemäß § 18 AO-SF durch die Entscheidung des Schulamtes **Wuppertal** der Bezirksregierung **Düsseldorf** vom **tt.mm.jjjj** den Förderschwerpunkt. &Er%Sie& wird zukünftig im Förderschwerpunkt **** gefördert.'
$ws.Range("C18").Value = 'FSP'
$ws.Range("B19").Value = '$Vorname$ wechselt gemäß § 17 AO-SF durch die Entscheidung des Schulamtes **Wuppertal** der Bezirksregierung **Düsseldorf** vom **tt.mm.jjjj** im Förderschwerpunkt *** den Bildungsgang. &Er%Sie& wird zukünftig im Bildungsgang **** gefördert.'
$ws.Range("C19").Value = 'FSP'
$ws.Range("B20").Value = 'Laut Beschluss der Klassenkonferenz vom **tt.mm.jjjj*** besteht gemäß § 17 AO-SF der Bedarf an sonderpädagogischer Unterstützung im Förderschwerpunkt **** mit dem zielgleichen Bildungsgang **** mit dem zieldifferenten **** weiterhin.'
$ws.Range("C20").Value = 'FSP'
$ws.Range("B21").Value = 'Ein Bericht zum Arbeits- und Sozialverhalten wird dem Zeugnis hinzugefügt, wenn die Versetzungskonferenz dies beschlossen hat und die Schulkonferenz dazu eine einheitliche Vorgehensweise festgelegt hat (§ 49 Absatz 2 Schulgesetzt NRW). Je nach Umfang kann dieser Bericht dem Zeugnis als Anlage hinzugefügt werden.'
$ws.Range("C21").Value = 'FSP'
$ws.Range("B22").Value = 'Wenn nach § 33 Absatz 3 AO-SF verfahren wird, werden die Noten in das Berichtszeugnis integriert. Es wird darauf hingewiesen, dass sich die Leistungsbewertung mit Noten an den Anforderungen der vorhergehenden Jahrgangsstufe der Grundschule oder der Hauptschule orientiert.'
$ws.Range("C22").Value = 'FSP'
$ws.Range("B23").Value = 'Die Zugehörigkeit zum Bildungsgang Lernen wurde gemäß § 18 AO-SF durch die Entscheidung des Schulamtes **Wuppertal** der Bezirksregierung **Düsseldorf** vom **tt.mm.jjjj** aufgehoben. Deshalb wird $Vorname§ zukünftig zielgleich im Bildungsgang der Grundschule unterrichtet. $Vorname$ hat aber weiterhin sonderpädagogischen Förderbedarf im Förderschwerpunkt ****.'
$ws.Range("C23").Value = 'FSP'
$ws.Range("B24").Value = '$Vorname$ nimmt im kommenden Schuljahr am Unterricht der Klasse 10 in einem besonderen Bildungsgang teil, mit dem Ziel, einem dem Ersten Schulabschluss gleichwertigen Abschluss zu erreichen.'
$ws.Range("C24").Value = 'FSP'
$ws.Range("B25").Value = '$Vorname$ hat den Abschluss des Bildungsgangs Lernen erworben.'
$ws.Range("C25").Value = 'FSP'
$ws.Range("B26").ClearContents()
$ws.Range("C26").Value = 'FSP'

# Column B ("Floskeltext") grows much wider to fit the longer remarks, and
# column C ("Art") gets its own (slightly wider) width now that it is no longer
# grouped together with column D.
$ws.Columns.Item(2).ColumnWidth = 254.5873
$ws.Columns.Item(3).ColumnWidth = 17.9225

# Leave the selection on the last edited row, matching the saved cursor position.
$ws.Range("B25").Select()
